$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-06-01"

# Update header label for the running-total "through" date column
$ws.Range("I1").Value = "2022 (through 06-01)"

# May value correction
$ws.Range("I5").Value = 115

# Add new June value
$ws.Range("I7").Value = 4

# Update yearly total
$ws.Range("I14").Value = 668
